$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# New "Sim 1" (column B) and "Sim2" (column C) values for rows 11-29.
# In the updated workbook both series end up with identical values.
$newValues = @(
    60.29,
    47.16,
    18.52,
    -8.94,
    12.39,
    -6.15,
    0.32,
    2.29,
    43.41,
    3.22,
    18.64,
    3.07,
    13.77,
    -2.29,
    13.13,
    -6.15,
    28.64,
    -6.66,
    -4.04
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = 11 + $i
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
    $ws.Cells.Item($row, 3).Value = $newValues[$i]
}

$ws.Range("D29").Select()
